# Book1.xlsx edit: add a new "Sheet6" with Find/Replace text-function examples,
# move the active selection, and leave Sheet2's selection parked on H3.

$wb = $excel.ActiveWorkbook

# Sheet2 was the previously-active sheet (tabSelected, selection A6). Move its
# selection to H3 before leaving it (it stops being the active/selected tab
# once Sheet6 is created & activated below).
$sheet2 = $wb.Worksheets.Item("Sheet2")
[void]$sheet2.Activate()
[void]$sheet2.Range("H3").Select()

# Add the new worksheet at the end of the tab strip (after the last existing
# sheet) and name it Sheet6.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Sheet6"

# Data + Find/Replace formula demos.
$ws.Range("A2").Value = "apple"
$ws.Range("A3").Value = "b"
$ws.Range("A4").Value = "c"
$ws.Range("C4").Formula = '=FIND("p",A2,1)'
$ws.Range("D4").Value = "Find"

$ws.Range("A5").Value = "d"
$ws.Range("C5").Formula = '=REPLACE(A2,2,1,"P")'
$ws.Range("D5").Value = "replace"

$ws.Range("A6").Value = "e"
$ws.Range("C6").Formula = '=REPLACE(FIND(A2,A2,1),1,1,"P")'

$ws.Range("C7").Formula = '=REPLACE(FIND("p",A2,1),1,2,"P")'

# Sheet6 becomes the active sheet/tab with selection on E6.
[void]$ws.Activate()
[void]$ws.Range("E6").Select()
